$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.2233361854305164
$ws.Range("E2").Value = 14.17398074038888
$ws.Range("F2").Value = 63.35973209733794
